# Update the "word" column (C) stimuli from the old placeholder words
# to the final German stimuli used for the first run of the experiment.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C2").Value = "de_folgen"
$ws.Range("C3").Value = "de_leeren"
$ws.Range("C4").Value = "de_warnen"

# Move the active selection from D6 to D7, matching the saved cursor
# position recorded in the sheet view.
$ws.Range("D7").Select()
